# Apply "Changes of 22nd June 2022" to Sheet1 of the workbook.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Update RouteWorkStartDate (P2) and FirstGenerationDate (CA2): 2022-06-17 -> 2022-06-22
$ws.Range("P2").Value = 44734
$ws.Range("CA2").Value = 44734

# Update RouteWorkReadyTime (S2) and FirstGenerationTime (CB2): 07:00 -> 22:00
$ws.Range("S2").Value = 0.91666666666666663

# Update RouteWorkScheduledEndTime (T2): 06:45 -> 21:45
$ws.Range("T2").Value = 0.90625

$ws.Range("CB2").Value = 0.91666666666666663

# Move the active selection to CD18, matching the saved view state in the workbook.
$ws.Activate()
$ws.Range("CD18").Select()
